# Applies the cryptos-list refresh described in the commit:
# "Updated cryptos list on Wed Apr 24 06:24:17 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.557.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.04%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.247.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.91%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.43%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.32%  "

# Row 7
$ws.Range("E7").Value = "  +0.10%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.244.28"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.86%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.546"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.52%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.24%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.91"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.86%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.505"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.17%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000270"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.86%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.86%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.784.33"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.01%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.630.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.11%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.37"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.38%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.257.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.23%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.114"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.12%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "505.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.44%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.72%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.749"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.64%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.25%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.82%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.97%  "

# Row 26
$ws.Range("E26").Value = "  +0.16%  "

# Row 27
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.28%  "

# Row 28
$ws.Range("B28").Value = "Hedera"
$ws.Range("C28").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.139"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +54.94%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.79%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.45%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.94%  "

# Row 32
$ws.Range("E32").Value = "  -8.18%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.49%  "

# Row 34
$ws.Range("E34").Value = "  +0.12%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.15"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.92%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.49%  "

# Row 37
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "55.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.97%  "

# Row 38
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0787"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +14.95%  "

# Row 39
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +17.72%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "492.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.15%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0426"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.98%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.129"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.68%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.13%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.292"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.88%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.36%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.946.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.11%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.85%  "

# Row 48
$ws.Range("E48").Value = "  +3.97%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.119"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.69%  "

# Row 50
$ws.Range("E50").Value = "  -0.03%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.23%  "
